$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the manager/associate email addresses from the old domain to the new one
for ($r = 2; $r -le 6; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    if ($bVal -ne $null -and $bVal -ne "") {
        $ws.Cells.Item($r, 2).Value2 = $bVal.Replace("@maveric-systems.com", "@gmail.com")
    }
    $dVal = $ws.Cells.Item($r, 4).Value2
    if ($dVal -ne $null -and $dVal -ne "") {
        $ws.Cells.Item($r, 4).Value2 = $dVal.Replace("@maveric-systems.com", "@gmail.com")
    }
}

# Update the active cell selection to D9 as recorded in the saved workbook view
$ws.Range("D9").Select()
